# Commit: "update due to suppressing oids"
#
# On the "Metadata" sheet:
#   - the "Experimental" property row (B7) now carries an explicit value
#     of the text "false" (it was previously left blank)
#   - the "Date" property row (B8) is refreshed to the new generation
#     timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Write B7 as literal text "false". A direct $ws.Range("B7").Value = "false"
# would be auto-coerced by Excel into a Boolean cell (t="b"), which is not
# what the source workbook stores (a shared-string "false"). Round-trip the
# literal through a formula-built string on a scratch cell, then paste only
# the resulting value back into B7 - this keeps the cell's existing style
# (s="2") and yields a genuine text cell.
$scratch = $ws.Cells.Item(1000, 50)
$scratch.Formula = '="false"'
$scratch.Copy()
$ws.Cells.Item(7, 2).PasteSpecial(-4163) # xlPasteValues
$scratch.Value = $null
$excel.CutCopyMode = 0

# Update the Date value.
$ws.Cells.Item(8, 2).Value = "2023-10-09T22:41:16+02:00"
